$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Save" column header in H1, using the same formatting as the
# existing header row (copy format only from G1 so the existing style is reused
# rather than a brand-new style being minted).
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null

# Add the data value for the new "Save" column
$ws.Range("H2").Value = 1
